$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (shifts old rows 12..36 down to 13..37).
$ws.Rows.Item(12).Insert()

# --- Row 11 ("Corrigeren (of geen)" / XYPositie ingemeten row) ---
# Pick up the correct cell borders/formatting by copying the style of row 10,
# which already uses the exact border pattern this row needs.
$ws.Range("A10:G10").Copy($ws.Range("A11:G11"))
$ws.Rows.Item(11).RowHeight = 30

$ws.Range("A11").Value = "Corrigeren (of geen)"
$ws.Range("B11").Value = "replaceRequest"
$ws.Range("C11").Value = "verplaatsverzoek"
$ws.Range("D11").Value = "verplaatsverzoek"
$ws.Range("E11").Value = "verplaatsverzoek"
$ws.Range("F11").Value = "XYPositie ingemeten"
$ws.Range("G11").Value = "Een nieuwe XY-inmeting hoeft niet noodzakelijkerwijs tot een correctie (verbeterde waarde) te leiden"

# --- Row 12 (new "Aanvullen en corrigeren" / XYZPositie ingemeten row) ---
# Copy the striped-row style (same pattern used by rows 4/6/13/...) onto the
# freshly inserted row.
$ws.Range("A4:G4").Copy($ws.Range("A12:G12"))

$ws.Range("A12").Value = "Aanvullen en corrigeren"
$ws.Range("B12").Value = "Combi"
$ws.Range("C12").Value = "Combi"
$ws.Range("D12").Value = "Combi"
$ws.Range("E12").Value = "Combi"
$ws.Range("F12").Value = "XYZPositie ingemeten"
$ws.Range("G12").Value = "Z en XY tegelijk ingemeten, combinatie van BRO-gebeurtenissen"

# Match the workbook's recorded selection after the edit.
$ws.Range("G12").Select()
